# product_dump.xlsx update:
#  - "filename" column values on the "product" sheet change from the
#    placeholder "temp" to "temp.png" (upload now carries a real extension)
#  - the worksheet's saved cell selection moves from H10 (just past the
#    used range) to H9 (the last populated cell)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("product")

# Column H ("filename") holds the shared string "temp" for every data row
# (rows 2-9). Rewrite each cell so the value becomes "temp.png".
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 8)
    if ($cell.Value2 -eq "temp") {
        $cell.Value = "temp.png"
    }
}

# Move/save the active selection on the sheet to H9.
$ws.Activate()
$ws.Range("H9").Select()
